# Update the "想去人数" (column F) figures for several events across the
# "展览" (Exhibition), "演出" (Performance) and "全部类型" (All types) sheets,
# matching the data refresh recorded in the commit "Update gh-pages to
# output generated at 456a3b4".
#
# NOTE: this engine's PowerShell dialect does not reliably support named
# (-Param value) arguments for user-defined functions, so everything below
# is called positionally.

$wb = $excel.ActiveWorkbook

function Set-ColumnFValues($SheetName, $RowValues) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $RowValues[$row]
    }
}

# Sheet "展览"
$exhibition = @{
    2  = 7641
    3  = 7641
    5  = 7826
    6  = 37
    8  = 28
    9  = 6554
    10 = 3348
    11 = 2
    12 = 3704
    14 = 42
    17 = 51
    20 = 9
    21 = 309
    23 = 3804
    27 = 280
    28 = 1452
    29 = 78
    30 = 51
    31 = 2730
    32 = 1767
    33 = 30
    34 = 41
    35 = 54
    36 = 3607
    37 = 287
    38 = 274
    41 = 527
    42 = 1398
    44 = 544
    45 = 632
}
Set-ColumnFValues "展览" $exhibition

# Sheet "演出"
$performance = @{
    13 = 87
    17 = 4
}
Set-ColumnFValues "演出" $performance

# Sheet "全部类型"
$allTypes = @{
    5  = 7641
    6  = 7641
    8  = 7826
    10 = 28
    11 = 6554
    12 = 3348
    13 = 2
    14 = 3704
    16 = 42
    19 = 51
    21 = 309
    24 = 3804
    31 = 280
    32 = 1452
    33 = 78
    34 = 51
    35 = 2730
    36 = 1767
    37 = 30
    38 = 41
    39 = 87
    40 = 3607
    41 = 287
    42 = 274
    45 = 527
    46 = 1398
    49 = 544
    50 = 632
}
Set-ColumnFValues "全部类型" $allTypes

$wb.Save()
